$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value2 = "Volume 29   Number  49"
$ws.Range("C9").Value2 = "Report Covering the Week  12/5/2022  Through  12/11/2022"

# --- Cells whose style/type changes (string <-> number) ---
$ws.Range("C16").Value2 = 2
$ws.Range("I28").Copy()
$ws.Range("C16").PasteSpecial(-4122)

$ws.Range("D16").Value2 = "'0"
$ws.Range("C23").Copy()
$ws.Range("D16").PasteSpecial(-4122)

$ws.Range("E16").Value2 = "'***.*"
$ws.Range("E23").Copy()
$ws.Range("E16").PasteSpecial(-4122)

$ws.Range("D17").Value2 = 1
$ws.Range("I28").Copy()
$ws.Range("D17").PasteSpecial(-4122)

$ws.Range("E17").Value2 = 100
$ws.Range("K28").Copy()
$ws.Range("E17").PasteSpecial(-4122)

$ws.Range("C18").Value2 = "'0"
$ws.Range("C23").Copy()
$ws.Range("C18").PasteSpecial(-4122)

$ws.Range("D27").Value2 = "'0"
$ws.Range("C23").Copy()
$ws.Range("D27").PasteSpecial(-4122)

$ws.Range("E27").Value2 = "'***.*"
$ws.Range("E23").Copy()
$ws.Range("E27").PasteSpecial(-4122)

$ws.Range("F27").Value2 = "'0"
$ws.Range("C23").Copy()
$ws.Range("F27").PasteSpecial(-4122)

# --- Cells whose style stays the same (value-only updates) ---
$ws.Range("F16").Value2 = 3
$ws.Range("G16").Value2 = 1
$ws.Range("H16").Value2 = 200
$ws.Range("I16").Value2 = 20
$ws.Range("K16").Value2 = 100
$ws.Range("L16").Value2 = 53.846153846153
$ws.Range("M16").Value2 = -9.090909090909
$ws.Range("N16").Value2 = -67.741935483871
$ws.Range("C17").Value2 = 2
$ws.Range("F17").Value2 = 6
$ws.Range("H17").Value2 = 500
$ws.Range("I17").Value2 = 49
$ws.Range("J17").Value2 = 31
$ws.Range("K17").Value2 = 58.064516129032
$ws.Range("L17").Value2 = -20.967741935483
$ws.Range("M17").Value2 = 4.255319148936
$ws.Range("N17").Value2 = -52.884615384615
$ws.Range("F18").Value2 = 7
$ws.Range("H18").Value2 = 600
$ws.Range("L18").Value2 = -4.651162790697
$ws.Range("M18").Value2 = -59.405940594059
$ws.Range("N18").Value2 = -87.1875
$ws.Range("C19").Value2 = 6
$ws.Range("D19").Value2 = 6
$ws.Range("E19").Value2 = 0
$ws.Range("F19").Value2 = 23
$ws.Range("G19").Value2 = 33
$ws.Range("H19").Value2 = -30.303030303030
$ws.Range("I19").Value2 = 256
$ws.Range("J19").Value2 = 185
$ws.Range("K19").Value2 = 38.378378378378
$ws.Range("L19").Value2 = 46.285714285714
$ws.Range("M19").Value2 = 80.281690140845
$ws.Range("N19").Value2 = 10.344827586206
$ws.Range("C20").Value2 = 1
$ws.Range("D20").Value2 = 1
$ws.Range("E20").Value2 = 0
$ws.Range("F20").Value2 = 11
$ws.Range("G20").Value2 = 7
$ws.Range("H20").Value2 = 57.142857142857
$ws.Range("I20").Value2 = 116
$ws.Range("J20").Value2 = 57
$ws.Range("K20").Value2 = 103.508771929825
$ws.Range("L20").Value2 = 163.636363636364
$ws.Range("M20").Value2 = 213.513513513514
$ws.Range("N20").Value2 = -83.016105417276
$ws.Range("C21").Value2 = 11
$ws.Range("D21").Value2 = 8
$ws.Range("E21").Value2 = 37.5
$ws.Range("F21").Value2 = 50
$ws.Range("G21").Value2 = 43
$ws.Range("H21").Value2 = 16.279069767441
$ws.Range("I21").Value2 = 486
$ws.Range("J21").Value2 = 308
$ws.Range("K21").Value2 = 57.792207792207
$ws.Range("L21").Value2 = 42.521994134897
$ws.Range("M21").Value2 = 37.677053824362
$ws.Range("N21").Value2 = -65.433854907539
$ws.Range("C24").Value2 = 9
$ws.Range("D24").Value2 = 7
$ws.Range("E24").Value2 = 28.571428571428
$ws.Range("F24").Value2 = 46
$ws.Range("G24").Value2 = 32
$ws.Range("H24").Value2 = 43.75
$ws.Range("I24").Value2 = 468
$ws.Range("J24").Value2 = 262
$ws.Range("K24").Value2 = 78.625954198473
$ws.Range("L24").Value2 = 73.977695167286
$ws.Range("M24").Value2 = -8.414872798434
$ws.Range("C25").Value2 = 2
$ws.Range("E25").Value2 = -50
$ws.Range("F25").Value2 = 13
$ws.Range("G25").Value2 = 18
$ws.Range("H25").Value2 = -27.777777777777
$ws.Range("I25").Value2 = 176
$ws.Range("J25").Value2 = 139
$ws.Range("K25").Value2 = 26.618705035971
$ws.Range("L25").Value2 = 35.384615384615
$ws.Range("M25").Value2 = -16.981132075471
$ws.Range("L26").Value2 = -60
$ws.Range("H27").Value2 = -100

$excel.CutCopyMode = 0
